# Truncate the "created_at" (column S) timestamps on the "All products"
# sheet down to plain dates, and force that column to be stored as text so
# Excel doesn't re-interpret the new "YYYY-MM-DD" strings as date serials.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All products")

# Make sure the whole created_at column (S, header + 6 data rows) is
# formatted as Text *before* writing the new values, otherwise Excel
# would silently turn "2021-05-15" into a date serial number again.
$ws.Range("S1:S7").NumberFormat = "@"

$ws.Range("S2").Value = "2021-05-15"
$ws.Range("S3").Value = "2021-05-15"
$ws.Range("S4").Value = "2021-05-15"
$ws.Range("S5").Value = "2021-05-15"
$ws.Range("S6").Value = "2021-05-15"
$ws.Range("S7").Value = "2021-05-16"

# Leave the cursor where it ended up after the last edit (row 8, col A)
# instead of the old selection on M7.
$ws.Range("A8").Select()
